# Dallas Away Pass Types - cleaning changes: insert a "Match ID" column at
# the front of the sheet (shifting every existing column one letter to the
# right) and stamp the match id (15) on every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column A; Excel shifts all existing columns (and their
# styles/number formats/shared-string references) one position to the
# right automatically.
$ws.Columns.Item(1).Insert()

# The header/blank rows (1-3) and the visible data rows (4-14) use a bold,
# borderless style for the new "Match ID" column - matches the workbook's
# existing header font (fontId 1) without the bordered/centered look used
# for the old "Player ID" column.
$ws.Range("A1:A14").Font.Bold = $true

# Header text for the new column.
$ws.Range("A1").Value2 = "Match ID"

# Every real data row (4 through the hidden totals row 15) belongs to the
# same match.
for ($r = 4; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value2 = 15
}

# Writing into the hidden rows nudges their row-height metadata; re-fit them
# back to the sheet default so the only persisted difference is the data.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(15).AutoFit()

# Leave the same kind of "whole visible data block" selection behind that
# was left after the real edit.
$ws.Range("A1:A14").Select()
